$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 221.71428
$ws.Range("I5").Value = 210.8421
$ws.Range("K5").Value = 210.8421
$ws.Range("M5").Value = -95.84209999999999
# Row 40
$ws.Range("H40").Value = 2678.35
$ws.Range("I40").Value = 2377.2
$ws.Range("J40").Value = 2778.7334
$ws.Range("K40").Value = 2377.2
$ws.Range("L40").Value = 2778.7334
$ws.Range("M40").Value = -2202.2
$ws.Range("N40").Value = -3128.7334
# Row 51
$ws.Range("H51").Value = 3429.7144
$ws.Range("J51").Value = 3401.6
$ws.Range("L51").Value = 3401.6
$ws.Range("N51").Value = -4369.6
# Row 103
$ws.Range("H103").Value = 920.75
$ws.Range("I103").Value = 561
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 1683
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -1097
$ws.Range("N103").Value = -7172
# Row 112
$ws.Range("H112").Value = 1725.75
$ws.Range("J112").Value = 1897.4
$ws.Range("L112").Value = 5692.200000000001
$ws.Range("N112").Value = -7908.200000000001
# Row 113
$ws.Range("H113").Value = 2999.75
$ws.Range("I113").Value = 2999.75
$ws.Range("K113").Value = 2999.75
$ws.Range("M113").Value = 254.25
# Row 118
$ws.Range("H118").Value = 2767.8
$ws.Range("I118").Value = 3280
$ws.Range("K118").Value = 9840
$ws.Range("M118").Value = -8183
# Row 132
$ws.Range("H132").Value = 43255.605
$ws.Range("I132").Value = 45377.527
$ws.Range("J132").Value = 5768.3335
$ws.Range("K132").Value = 136132.581
$ws.Range("L132").Value = 17305.0005
$ws.Range("M132").Value = -133602.581
$ws.Range("N132").Value = -22365.0005
# Row 137
$ws.Range("H137").Value = 290838.22
$ws.Range("I137").Value = 397937.66
$ws.Range("J137").Value = 2493.6155
$ws.Range("K137").Value = 1193812.98
$ws.Range("L137").Value = 7480.8465
$ws.Range("M137").Value = -1191262.98
$ws.Range("N137").Value = -12580.8465
# Row 138
$ws.Range("H138").Value = 3117.6428
$ws.Range("J138").Value = 4324.8184
$ws.Range("L138").Value = 12974.4552
$ws.Range("N138").Value = -23254.4552
# Row 141
$ws.Range("H141").Value = 1437.1666
$ws.Range("I141").Value = 1437.1666
$ws.Range("K141").Value = 4311.4998
$ws.Range("M141").Value = 868.5002000000004

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3668151.5
$ws.Range("I32").Value = 4507296
$ws.Range("K32").Value = 4507296
$ws.Range("M32").Value = -4507009
# Row 45
$ws.Range("H45").Value = 4223.864
$ws.Range("I45").Value = 4557.3335
$ws.Range("K45").Value = 4557.3335
$ws.Range("M45").Value = -4180.3335
# Row 61
$ws.Range("H61").Value = 1765326.8
$ws.Range("I61").Value = 1765326.8
$ws.Range("K61").Value = 1765326.8
$ws.Range("M61").Value = -1765114.8
# Row 74
$ws.Range("H74").Value = 2317582.5
$ws.Range("I74").Value = 3380639
$ws.Range("K74").Value = 3380639
$ws.Range("M74").Value = -3379765
# Row 77
$ws.Range("H77").Value = 2317582.5
$ws.Range("I77").Value = 3380639
$ws.Range("K77").Value = 16903195
$ws.Range("M77").Value = -16898827
# Row 92
$ws.Range("H92").Value = 55477.09
$ws.Range("J92").Value = 55477.09
$ws.Range("L92").Value = 55477.09
$ws.Range("N92").Value = -60469.09
# Row 110
$ws.Range("H110").Value = 740.1
$ws.Range("I110").Value = 655.6667
$ws.Range("K110").Value = 655.6667
$ws.Range("M110").Value = 1389.3333
# Row 122
$ws.Range("H122").Value = 3355
$ws.Range("I122").Value = 1853
$ws.Range("J122").Value = 4857
$ws.Range("K122").Value = 5559
$ws.Range("L122").Value = 14571
$ws.Range("M122").Value = -3109
$ws.Range("N122").Value = -19471
# Row 132
$ws.Range("H132").Value = 600016.9
$ws.Range("I132").Value = 707395.8
$ws.Range("K132").Value = 2122187.4
$ws.Range("M132").Value = -2119657.4
# Row 136
$ws.Range("H136").Value = 1765326.8
$ws.Range("I136").Value = 1765326.8
$ws.Range("K136").Value = 5295980.4
$ws.Range("M136").Value = -5293430.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2666.6667
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
# Row 89
$ws.Range("H89").Value = 2666.6667
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
# Row 94
$ws.Range("H94").Value = 2049.8
$ws.Range("I94").Value = 1166.5555
$ws.Range("K94").Value = 1166.5555
$ws.Range("M94").Value = -715.5554999999999
# Row 99
$ws.Range("H99").Value = 130247
$ws.Range("I99").Value = 171329.67
$ws.Range("K99").Value = 171329.67
$ws.Range("M99").Value = -169831.67
# Row 107
$ws.Range("H107").Value = 1594.4688
$ws.Range("I107").Value = 1551.1072
$ws.Range("K107").Value = 1551.1072
$ws.Range("M107").Value = 368.8928000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7257.5625
$ws.Range("I31").Value = 2733.4119
$ws.Range("K31").Value = 2733.4119
$ws.Range("M31").Value = -2438.4119
# Row 34
$ws.Range("H34").Value = 7257.5625
$ws.Range("I34").Value = 2733.4119
$ws.Range("K34").Value = 2733.4119
$ws.Range("M34").Value = -2531.4119
# Row 58
$ws.Range("H58").Value = 1375259.4
$ws.Range("I58").Value = 1765476.2
$ws.Range("K58").Value = 1765476.2
$ws.Range("M58").Value = -1765273.2
# Row 136
$ws.Range("H136").Value = 1375259.4
$ws.Range("I136").Value = 1765476.2
$ws.Range("K136").Value = 5296428.6
$ws.Range("M136").Value = -5293878.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2978.7646
$ws.Range("J5").Value = 7123.8335
$ws.Range("L5").Value = 21371.5005
$ws.Range("N5").Value = -21595.5005
# Row 80
$ws.Range("H80").Value = 5328
$ws.Range("I80").Value = 4994.5
$ws.Range("K80").Value = 14983.5
$ws.Range("M80").Value = -14047.5
# Row 83
$ws.Range("H83").Value = 5328
$ws.Range("I83").Value = 4994.5
$ws.Range("K83").Value = 44950.5
$ws.Range("M83").Value = -40270.5
# Row 131
$ws.Range("H131").Value = 7437.15
$ws.Range("I131").Value = 1484
$ws.Range("K131").Value = 4452
$ws.Range("M131").Value = 588
# Row 135
$ws.Range("H135").Value = 2978.7646
$ws.Range("J135").Value = 7123.8335
$ws.Range("L135").Value = 64114.5015
$ws.Range("N135").Value = -69184.5015

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 215.5
$ws.Range("I2").Value = 114.625
$ws.Range("J2").Value = 273.14285
$ws.Range("K2").Value = 114.625
$ws.Range("L2").Value = 273.14285
$ws.Range("M2").Value = -1.625
$ws.Range("N2").Value = -499.14285
# Row 70
$ws.Range("H70").Value = 7915.5
$ws.Range("I70").Value = 8672.5
$ws.Range("K70").Value = 8672.5
$ws.Range("M70").Value = -8402.5
# Row 73
$ws.Range("H73").Value = 7915.5
$ws.Range("I73").Value = 8672.5
$ws.Range("K73").Value = 8672.5
$ws.Range("M73").Value = -7736.5
# Row 102
$ws.Range("H102").Value = 2048.96
$ws.Range("I102").Value = 1486.8572
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1486.8572
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 135.1428000000001
$ws.Range("N102").Value = -8244
# Row 107
$ws.Range("H107").Value = 254
$ws.Range("I107").Value = 165
$ws.Range("J107").Value = 343
$ws.Range("K107").Value = 165
$ws.Range("L107").Value = 343
$ws.Range("M107").Value = 1755
$ws.Range("N107").Value = -4183
# Row 122
$ws.Range("H122").Value = 11800
$ws.Range("I122").Value = 11800
$ws.Range("K122").Value = 35400
$ws.Range("M122").Value = -32950
# Row 132
$ws.Range("H132").Value = 2360.9167
$ws.Range("I132").Value = 1484.625
$ws.Range("J132").Value = 4113.5
$ws.Range("K132").Value = 4453.875
$ws.Range("L132").Value = 12340.5
$ws.Range("M132").Value = -1923.875
$ws.Range("N132").Value = -17400.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 210
$ws.Range("I55").Value = 173.11111
$ws.Range("J55").Value = 257.42856
$ws.Range("K55").Value = 173.11111
$ws.Range("L55").Value = 257.42856
$ws.Range("M55").Value = -0.1111099999999965
$ws.Range("N55").Value = -603.4285600000001
# Row 100
$ws.Range("H100").Value = 9890.214
$ws.Range("I100").Value = 2933.25
$ws.Range("J100").Value = 19166.166
$ws.Range("K100").Value = 2933.25
$ws.Range("L100").Value = 19166.166
$ws.Range("M100").Value = -2392.25
$ws.Range("N100").Value = -20248.166
# Row 132
$ws.Range("H132").Value = 3792
$ws.Range("I132").Value = 3304
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 9912
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -7382
$ws.Range("N132").Value = -21560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 56245.715
$ws.Range("J95").Value = 56245.715
$ws.Range("L95").Value = 56245.715
$ws.Range("N95").Value = -61737.715
# Row 107
$ws.Range("H107").Value = 2847.9
$ws.Range("I107").Value = 1337.1
$ws.Range("K107").Value = 4011.3
$ws.Range("M107").Value = -2091.3
# Row 126
$ws.Range("H126").Value = 4403.1763
$ws.Range("I126").Value = 3681.5386
$ws.Range("K126").Value = 11044.6158
$ws.Range("M126").Value = -8574.6158
# Row 132
$ws.Range("H132").Value = 3971991
$ws.Range("I132").Value = 5211391
$ws.Range("K132").Value = 15634173
$ws.Range("M132").Value = -15631643
